$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the old "total" row (row 63): drop A63/C63 entirely, and
#        blank out B63's formula+style so it becomes a bare empty cell. ---
$ws.Range("A63").Clear()
$ws.Range("C63").Clear()
$ws.Range("B63").ClearContents()
$ws.Range("B63").Style = "Normal"

# --- 2. Rows 56/57: these rows already hold data, only their A/C cell
#        formatting changes (duplicate "Normal"-style slot). ---
$ws.Range("A56").Style = "Normal"
$ws.Range("C56").Style = "Normal"
$ws.Range("A57").Style = "Normal"
$ws.Range("C57").Style = "Normal"

# --- 3. Fill in the new Docker-troubleshooting entries in rows 58-61.
#        Date-like strings are forced to text with a leading apostrophe so
#        they aren't silently converted into date serials, then the style
#        is normalized back to plain (no quote-prefix) formatting. ---
$ws.Range("A58").Value = "'03.01.2025"
$ws.Range("A58").Style = "Normal"
$ws.Range("B58").Value = 4
$ws.Range("C58").Value = "Docker troubleshooting (lief ember mismatch)"

$ws.Range("B59").Value = 4
$ws.Range("C59").Value = "Ember neu trainieren"

$ws.Range("B60").Value = 2
$ws.Range("C60").Value = "Troubleshooting Lief"

$ws.Range("A61").Value = "'04.01.2025"
$ws.Range("A61").Style = "Normal"
$ws.Range("B61").Value = 3

# --- 4. Add the new total row at row 72 (SUM now covers B2:B71), reusing
#        the existing "Check Cell" style used by the old total row. ---
$ws.Range("A72").Value = "total"
$ws.Range("A72").Style = "Check Cell"
$ws.Range("B72").Formula = "=SUM(B2:B71)"
$ws.Range("B72").Style = "Check Cell"

# --- 5. Reset the saved scroll position back to the top of the sheet. ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
